$p = $ppt.ActivePresentation

function Find-ShapeWithText($slide, [string]$needle) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text.IndexOf($needle) -ge 0) {
                return $sh
            }
        }
    }
    return $null
}

# --- Slide 16: "Our metrics could be baised by the task choice!"
#     -> "Our metrics could be biased by the task choice!"
#     (spelling fix; PowerPoint splits the corrected word into its own run) ---
$s16 = $p.Slides.Item(16)
$shape16 = Find-ShapeWithText $s16 "baised"
$tr16 = $shape16.TextFrame.TextRange
$full16 = $tr16.Text
$idx16 = $full16.IndexOf("baised")
$word16 = $tr16.Characters($idx16 + 1, 6)
$word16.Text = "biased"

# --- Slide 8: "For each driver predict the if they will DNF (did not finish) a race in the next 1 month"
#     -> "For each driver predict if they will DNF (did not finish) a race in the next 1 month"
#     (drop the stray "the") ---
$s8 = $p.Slides.Item(8)
$shape8 = Find-ShapeWithText $s8 "DNF"
$tr8 = $shape8.TextFrame.TextRange
$full8 = $tr8.Text
$runStart8 = $full8.IndexOf("For each driver") + 1
$runLen8 = $full8.Length - $runStart8 + 1
$run8 = $tr8.Characters($runStart8, $runLen8)
$run8.Text = "For each driver predict if they will DNF (did not finish) a race in the next 1 month"
